$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value2 = 370.13635
$ws.Range("I19").Value2 = 256.75
$ws.Range("K19").Value2 = 256.75
$ws.Range("M19").Value2 = -81.75
# Row 40
$ws.Range("H40").Value2 = 3599.0908
$ws.Range("I40").Value2 = 3240.4707
$ws.Range("J40").Value2 = 4818.4
$ws.Range("K40").Value2 = 3240.4707
$ws.Range("L40").Value2 = 4818.4
$ws.Range("M40").Value2 = -3065.4707
$ws.Range("N40").Value2 = -5168.4
# Row 88
$ws.Range("H88").Value2 = 892.5
$ws.Range("I88").Value2 = 900
$ws.Range("J88").Value2 = 890
$ws.Range("K88").Value2 = 900
$ws.Range("L88").Value2 = 890
$ws.Range("M88").Value2 = -494
$ws.Range("N88").Value2 = -1702
# Row 91
$ws.Range("H91").Value2 = 892.5
$ws.Range("I91").Value2 = 900
$ws.Range("J91").Value2 = 890
$ws.Range("K91").Value2 = 900
$ws.Range("L91").Value2 = 890
$ws.Range("M91").Value2 = 504
$ws.Range("N91").Value2 = -3698
# Row 113
$ws.Range("H113").Value2 = 16332.667
$ws.Range("J113").Value2 = 12749.25
$ws.Range("L113").Value2 = 12749.25
$ws.Range("N113").Value2 = -19257.25
# Row 114
$ws.Range("H114").Value2 = 89500
$ws.Range("J114").Value2 = 89500
$ws.Range("L114").Value2 = 89500
$ws.Range("N114").Value2 = -98178
# Row 138
$ws.Range("H138").Value2 = 2313.4167
$ws.Range("I138").Value2 = 822.2258
$ws.Range("J138").Value2 = 3440.9023
$ws.Range("K138").Value2 = 2466.6774
$ws.Range("L138").Value2 = 10322.7069
$ws.Range("M138").Value2 = 2673.3226
$ws.Range("N138").Value2 = -20602.7069

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 5302
$ws.Range("I32").Value2 = 5233.913
$ws.Range("K32").Value2 = 5233.913
$ws.Range("M32").Value2 = -4946.913
# Row 61
$ws.Range("H61").Value2 = 13790.6
$ws.Range("I61").Value2 = 17460.363
$ws.Range("J61").Value2 = 3698.75
$ws.Range("K61").Value2 = 17460.363
$ws.Range("L61").Value2 = 3698.75
$ws.Range("M61").Value2 = -17248.363
$ws.Range("N61").Value2 = -4122.75
# Row 97
$ws.Range("H97").Value2 = 8004590
$ws.Range("I97").Value2 = 7163.1333
$ws.Range("J97").Value2 = 20000730
$ws.Range("K97").Value2 = 7163.1333
$ws.Range("L97").Value2 = 20000730
$ws.Range("M97").Value2 = -6667.1333
$ws.Range("N97").Value2 = -20001722
# Row 122
$ws.Range("H122").Value2 = 1298670.9
$ws.Range("I122").Value2 = 4172.36
$ws.Range("K122").Value2 = 12517.08
$ws.Range("M122").Value2 = -10067.08
# Row 132
$ws.Range("H132").Value2 = 3640.28
$ws.Range("I132").Value2 = 3014.375
$ws.Range("K132").Value2 = 9043.125
$ws.Range("M132").Value2 = -6513.125
# Row 136
$ws.Range("H136").Value2 = 13790.6
$ws.Range("I136").Value2 = 17460.363
$ws.Range("J136").Value2 = 3698.75
$ws.Range("K136").Value2 = 52381.08900000001
$ws.Range("L136").Value2 = 11096.25
$ws.Range("M136").Value2 = -49831.08900000001
$ws.Range("N136").Value2 = -16196.25
# Row 139
$ws.Range("H139").Value2 = 191005
$ws.Range("J139").Value2 = 191005
$ws.Range("L139").Value2 = 191005
$ws.Range("N139").Value2 = -201285

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 98
$ws.Range("H98").Value2 = 542000
$ws.Range("J98").Value2 = 542000
$ws.Range("L98").Value2 = 542000
$ws.Range("N98").Value2 = -547990
# Row 99
$ws.Range("H99").Value2 = 17341.44
$ws.Range("I99").Value2 = 21207.277
$ws.Range("K99").Value2 = 21207.277
$ws.Range("M99").Value2 = -19709.277
# Row 105
$ws.Range("H105").Value2 = 104523
$ws.Range("I105").Value2 = 200646
$ws.Range("J105").Value2 = 8400
$ws.Range("K105").Value2 = 200646
$ws.Range("L105").Value2 = 8400
$ws.Range("M105").Value2 = -198899
$ws.Range("N105").Value2 = -11894
# Row 134
$ws.Range("H134").Value2 = 8771.764999999999
$ws.Range("I134").Value2 = 9774.733
$ws.Range("J134").Value2 = 1249.5
$ws.Range("K134").Value2 = 29324.199
$ws.Range("L134").Value2 = 3748.5
$ws.Range("M134").Value2 = -26789.199
$ws.Range("N134").Value2 = -8818.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 16093.5
$ws.Range("I31").Value2 = 18366.875
$ws.Range("J31").Value2 = 7000
$ws.Range("K31").Value2 = 18366.875
$ws.Range("L31").Value2 = 7000
$ws.Range("M31").Value2 = -18071.875
$ws.Range("N31").Value2 = -7590
# Row 34
$ws.Range("H34").Value2 = 16093.5
$ws.Range("I34").Value2 = 18366.875
$ws.Range("J34").Value2 = 7000
$ws.Range("K34").Value2 = 18366.875
$ws.Range("L34").Value2 = 7000
$ws.Range("M34").Value2 = -18164.875
$ws.Range("N34").Value2 = -7404
# Row 99
$ws.Range("H99").Value2 = 11367425
$ws.Range("I99").Value2 = 13891742
$ws.Range("K99").Value2 = 13891742
$ws.Range("M99").Value2 = -13890244
# Row 105
$ws.Range("H105").Value2 = 236260.33
$ws.Range("I105").Value2 = 301334.72
$ws.Range("K105").Value2 = 301334.72
$ws.Range("M105").Value2 = -299587.72
# Row 114
$ws.Range("H114").Value2 = 11555
$ws.Range("J114").Value2 = 11555
$ws.Range("L114").Value2 = 11555
$ws.Range("N114").Value2 = -20233
# Row 126
$ws.Range("H126").Value2 = 11367425
$ws.Range("I126").Value2 = 13891742
$ws.Range("K126").Value2 = 41675226
$ws.Range("M126").Value2 = -41672756

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 35
$ws.Range("H35").Value2 = 982.2222
$ws.Range("J35").Value2 = 1133.3334
$ws.Range("L35").Value2 = 3400.0002
$ws.Range("N35").Value2 = -3976.0002
# Row 97
$ws.Range("H97").Value2 = 200100
$ws.Range("I97").Value2 = 300000
$ws.Range("K97").Value2 = 900000
$ws.Range("M97").Value2 = -899504
# Row 103
$ws.Range("H103").Value2 = 6574.1177
$ws.Range("I103").Value2 = 8974.833000000001
$ws.Range("J103").Value2 = 5264.636
$ws.Range("K103").Value2 = 26924.499
$ws.Range("L103").Value2 = 15793.908
$ws.Range("M103").Value2 = -26045.499
$ws.Range("N103").Value2 = -17551.908
# Row 107
$ws.Range("H107").Value2 = 1040.2693
$ws.Range("I107").Value2 = 385.7143
$ws.Range("J107").Value2 = 1281.421
$ws.Range("K107").Value2 = 1157.1429
$ws.Range("L107").Value2 = 3844.263
$ws.Range("M107").Value2 = 762.8571000000002
$ws.Range("N107").Value2 = -7684.263
# Row 115
$ws.Range("H115").Value2 = 1956.3334
$ws.Range("J115").Value2 = 4974
$ws.Range("L115").Value2 = 14922
$ws.Range("N115").Value2 = -17272
# Row 131
$ws.Range("H131").Value2 = 1613.5469
$ws.Range("I131").Value2 = 1299
$ws.Range("J131").Value2 = 1652.1754
$ws.Range("K131").Value2 = 3897
$ws.Range("L131").Value2 = 4956.5262
$ws.Range("M131").Value2 = 1143
$ws.Range("N131").Value2 = -15036.5262

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value2 = 106.37037
$ws.Range("I2").Value2 = 80.166664
$ws.Range("J2").Value2 = 316
$ws.Range("K2").Value2 = 80.166664
$ws.Range("L2").Value2 = 316
$ws.Range("M2").Value2 = 32.833336
$ws.Range("N2").Value2 = -542
# Row 132
$ws.Range("H132").Value2 = 4507.364
$ws.Range("I132").Value2 = 4458.1
$ws.Range("K132").Value2 = 13374.3
$ws.Range("M132").Value2 = -10844.3

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value2 = 6957.8125
$ws.Range("I22").Value2 = 11490.111
$ws.Range("K22").Value2 = 11490.111
$ws.Range("M22").Value2 = -11195.111
# Row 27
$ws.Range("H27").Value2 = 6957.8125
$ws.Range("I27").Value2 = 11490.111
$ws.Range("K27").Value2 = 11490.111
$ws.Range("M27").Value2 = -11383.111
# Row 68
$ws.Range("H68").Value2 = 6865.8335
$ws.Range("I68").Value2 = 1696
$ws.Range("K68").Value2 = 1696
$ws.Range("M68").Value2 = -947
# Row 71
$ws.Range("H71").Value2 = 6865.8335
$ws.Range("I71").Value2 = 1696
$ws.Range("K71").Value2 = 8480
$ws.Range("M71").Value2 = -4736
# Row 100
$ws.Range("H100").Value2 = 7499.4
$ws.Range("I100").Value2 = 2999.6667
$ws.Range("K100").Value2 = 2999.6667
$ws.Range("M100").Value2 = -2458.6667
# Row 122
$ws.Range("H122").Value2 = 7493.6
$ws.Range("I122").Value2 = 11985
$ws.Range("J122").Value2 = 4499.3335
$ws.Range("K122").Value2 = 35955
$ws.Range("L122").Value2 = 13498.0005
$ws.Range("M122").Value2 = -33505
$ws.Range("N122").Value2 = -18398.0005
# Row 132
$ws.Range("H132").Value2 = 500044.1
$ws.Range("I132").Value2 = 623867.7
$ws.Range("K132").Value2 = 1871603.1
$ws.Range("M132").Value2 = -1869073.1

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value2 = 13600
$ws.Range("I45").Value2 = 0
$ws.Range("J45").Value2 = 13600
$ws.Range("K45").Value2 = 0
$ws.Range("L45").Value2 = 13600
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value2 = -14582
# Row 107
$ws.Range("H107").Value2 = 19350.295
$ws.Range("I107").Value2 = 1710.8572
$ws.Range("K107").Value2 = 5132.571599999999
$ws.Range("M107").Value2 = -3212.571599999999
# Row 122
$ws.Range("H122").Value2 = 5081.8213
$ws.Range("I122").Value2 = 1814.6
$ws.Range("K122").Value2 = 5443.799999999999
$ws.Range("M122").Value2 = -2993.799999999999
# Row 132
$ws.Range("H132").Value2 = 10279.631
$ws.Range("I132").Value2 = 11082.171
$ws.Range("J132").Value2 = 3698.8
$ws.Range("K132").Value2 = 33246.513
$ws.Range("L132").Value2 = 11096.4
$ws.Range("M132").Value2 = -30716.513
$ws.Range("N132").Value2 = -16156.4
# Row 141
$ws.Range("H141").Value2 = 111051.445
$ws.Range("J141").Value2 = 106307.875
$ws.Range("L141").Value2 = 106307.875
$ws.Range("N141").Value2 = -116667.875
